$d = $word.ActiveDocument
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$rng = $lastPara.Range
$rng.Collapse(1)
$rng.InsertParagraphBefore()
$target = $d.Paragraphs.Item($n)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Refrences </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve"> Alahmari, M., &amp; Khalil, I. (2021). Building scalable eCommerce web applications with Next.js and React. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>International Journal of Computer Science and Network Security</w:t></w:r><w:r><w:t>, 21(4), 45–52. https://doi.org/10.22937/IJCSNS.2021.21.4.6</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Banks, A. (2020). </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>React Native in Action</w:t></w:r><w:r><w:t xml:space="preserve"> (1st ed.). Manning Publications.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Bierman, G., Abadi, M., &amp; Torgersen, M. (2014). Understanding TypeScript. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Proceedings of the ACM on Programming Languages</w:t></w:r><w:r><w:t>, 1(ICFP), 1–20. https://doi.org/10.1145/2628136</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Chen, J., &amp; Lin, Y. (2022). Leveraging Next.js for server-side rendering in eCommerce platforms. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Journal of Web Engineering</w:t></w:r><w:r><w:t>, 21(2), 115–132. https://doi.org/10.1145/3456789</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Gravina, D., &amp; Mastroeni, L. (2021). Building cross-platform mobile apps with React Native and TypeScript. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Software: Practice and Experience</w:t></w:r><w:r><w:t>, 51(12), 2665–2681. https://doi.org/10.1002/spe.2961</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Hossain, M., &amp; Islam, M. (2020). A study on UI/UX design principles for eCommerce mobile applications. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>International Journal of Computer Applications</w:t></w:r><w:r><w:t>, 175(9), 15–22. https://doi.org/10.5120/ijca2020919677</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Kumar, R., &amp; Singh, A. (2021). Secure payment integration strategies in modern eCommerce apps. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>International Journal of Advanced Research in Computer Science</w:t></w:r><w:r><w:t>, 12(5), 65–73. https://doi.org/10.26483/ijarcs.v12i5.7043</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Li, Q., &amp; Wang, Y. (2021). Real-time inventory management in online retail systems. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Journal of Retailing and Consumer Services</w:t></w:r><w:r><w:t>, 58, 102300. https://doi.org/10.1016/j.jretconser.2020.102300</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Rajput, N., &amp; Patil, S. (2022). Implementing product recommendation engines in eCommerce platforms using React. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>International Journal of Engineering Research &amp; Technology</w:t></w:r><w:r><w:t>, 11(3), 144–150. https://doi.org/10.17577/IJERTV11IS030218</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Singh, P., &amp; Sharma, R. (2020). Mobile-first approach for eCommerce application development using React Native. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>International Journal of Computer Science Trends and Technology</w:t></w:r><w:r><w:t>, 8(1), 50–58.</w:t></w:r></w:p>'
$target.Range.InsertXML($xml)
Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
